# Update "想去人数" (want-to-go count) values in column F across the
# 展览 (Exhibition), 演出 (Performance) and 全部类型 (All types) sheets.
# 本地生活 (Local life) sheet is unchanged.

$wb = $excel.ActiveWorkbook

# --- 展览 sheet -------------------------------------------------------
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value  = 198
$ws1.Range("F5").Value  = 975
$ws1.Range("F6").Value  = 5391
$ws1.Range("F7").Value  = 469
$ws1.Range("F8").Value  = 657
$ws1.Range("F9").Value  = 933
$ws1.Range("F11").Value = 74
$ws1.Range("F12").Value = 34
$ws1.Range("F13").Value = 579
$ws1.Range("F17").Value = 1795
$ws1.Range("F18").Value = 1460
$ws1.Range("F19").Value = 867
$ws1.Range("F22").Value = 320
$ws1.Range("F23").Value = 530
$ws1.Range("F24").Value = 141
$ws1.Range("F25").Value = 1051
$ws1.Range("F28").Value = 2736
$ws1.Range("F31").Value = 60
$ws1.Range("F32").Value = 107
$ws1.Range("F34").Value = 343
$ws1.Range("F38").Value = 227
$ws1.Range("F39").Value = 282
$ws1.Range("F40").Value = 671
$ws1.Range("F43").Value = 52

# --- 演出 sheet -------------------------------------------------------
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F4").Value = 167
$ws2.Range("F6").Value = 117

# --- 全部类型 sheet ----------------------------------------------------
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F3").Value  = 198
$ws4.Range("F5").Value  = 975
$ws4.Range("F7").Value  = 5391
$ws4.Range("F8").Value  = 469
$ws4.Range("F9").Value  = 657
$ws4.Range("F11").Value = 167
$ws4.Range("F12").Value = 933
$ws4.Range("F15").Value = 117
$ws4.Range("F16").Value = 74
$ws4.Range("F17").Value = 34
$ws4.Range("F18").Value = 579
$ws4.Range("F23").Value = 1795
$ws4.Range("F24").Value = 1460
$ws4.Range("F25").Value = 867
$ws4.Range("F27").Value = 320
$ws4.Range("F29").Value = 530
$ws4.Range("F30").Value = 141
$ws4.Range("F31").Value = 1051
$ws4.Range("F32").Value = 2737
$ws4.Range("F35").Value = 60
$ws4.Range("F36").Value = 107
$ws4.Range("F38").Value = 343
$ws4.Range("F42").Value = 282
$ws4.Range("F43").Value = 671
$ws4.Range("F45").Value = 52

$wb.Save()
